# Applies the "Add files via upload" / "Data files updated June 8, 2020 19:05"
# revision to the COVID case-tracking workbook:
#   - updates the "DHEC Daily Cases" (column E) figures for a batch of
#     already-recorded dates (the source agency revised several daily counts)
#   - fills in the previously-missing E/F figures for 6/7/2020 (row 95)
#   - appends the new day, 6/8/2020, as row 96
#   - extends the two chart series on the "Chart" sheet so they keep
#     covering the full data range
#   - leaves the "Data" sheet as the active/selected tab (it was "Chart")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------
# 1) Revised "DHEC Daily Cases" (column E) values for existing rows.
#    Column F ("DHEC Sum") is a running total formula (=SUM(prevF,E)) so
#    it -- and every downstream dependent -- recalculates automatically.
# ---------------------------------------------------------------------
$eUpdates = @{
    24 = 111
    27 = 205
    37 = 139
    39 = 121
    41 = 89
    42 = 266
    52 = 150
    54 = 108
    56 = 198
    57 = 156
    58 = 200
    59 = 136
    64 = 232
    66 = 99
    69 = 130
    70 = 169
    72 = 276
    73 = 156
    76 = 141
    78 = 250
    79 = 232
    80 = 190
    81 = 90
    82 = 255
    83 = 205
    85 = 323
    86 = 404
    87 = 312
    88 = 269
    89 = 262
    90 = 226
    91 = 353
    92 = 423
    93 = 507
    94 = 382
}

foreach ($r in $eUpdates.Keys) {
    $ws.Cells.Item($r, 5).Value = $eUpdates[$r]
}

# ---------------------------------------------------------------------
# 2) Row 95 (6/7/2020) previously had no DHEC daily figure yet -- fill in
#    E95 and the matching running-total formula in F95.
# ---------------------------------------------------------------------
$ws.Range("E95").Value = 542
$ws.Range("F95").Formula = "=SUM(F94,E95)"

# ---------------------------------------------------------------------
# 3) Append the new day: row 96 (6/8/2020).
# ---------------------------------------------------------------------
$ws.Range("A96").Value = 43990
$ws.Range("A95").Copy()
$ws.Range("A96").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B96").Value = 14800
$ws.Range("C96").Value = 557
$ws.Range("D96").Formula = "=SUM(C96,-(C95))"
$ws.Range("H96").Value = 23401
$ws.Range("I96").Value = 229861

$ws.Range("J95").Copy()
$ws.Range("J96").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("J96").Formula = "=IMDIV(H96,I96)"

# ---------------------------------------------------------------------
# 4) Extend the two "Chart" sheet series (and their shared date
#    category range) from row 95 out to the new row 96.
# ---------------------------------------------------------------------
$chartSheet = $wb.Worksheets.Item("Chart")
$chart = $chartSheet.ChartObjects(1).Chart

$series1 = $chart.SeriesCollection(1)
$series1.Formula = "=SERIES(Data!`$B`$1,Data!`$A`$2:`$A`$96,Data!`$B`$2:`$B`$96,1)"

$series2 = $chart.SeriesCollection(2)
$series2.Formula = "=SERIES(Data!`$C`$1,Data!`$A`$2:`$A`$96,Data!`$C`$2:`$C`$96,2)"

# ---------------------------------------------------------------------
# 5) Leave the "Data" sheet active/selected (previously "Chart" was the
#    selected tab) with the view scrolled down near the new rows.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K74").Select()
